# Update column G ("K") values on Sheet1 (rows 2-29) to reflect the
# regenerated save_data values (K instead of Strike#).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$kValues = @{
    2  = 1
    3  = 2
    4  = 1
    5  = 3
    6  = 1
    7  = 3
    8  = 2
    9  = 1
    10 = 0
    11 = 2
    12 = 1
    13 = 0
    14 = 0
    15 = 2
    16 = 0
    17 = 3
    18 = 2
    19 = 0
    20 = 2
    21 = 1
    22 = 0
    23 = 0
    24 = 2
    25 = 5
    26 = 4
    27 = 0
    28 = 3
    29 = 0
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
